$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.621.20'
$ws.Range('E2').Value = '  +0.59%  '

$ws.Range('D3').Value = '1.636.36'
$ws.Range('E3').Value = '  -0.49%  '

$ws.Range('E4').Value = '  +0.04%  '

$ws.Range('D5').Value = '212.39'
$ws.Range('E5').Value = '  -0.01%  '

$ws.Range('D6').Value = '0.523'
$ws.Range('E6').Value = '  -1.02%  '

$ws.Range('E7').Value = '  +0.04%  '

$ws.Range('D8').Value = '22.92'
$ws.Range('E8').Value = '  -1.82%  '

$ws.Range('E9').Value = '  +0.02%  '

$ws.Range('E10').Value = '  -0.14%  '

$ws.Range('D11').Value = '0.0892'
$ws.Range('E11').Value = '  +0.09%  '

$ws.Range('D12').Value = '1.868.50'
$ws.Range('E12').Value = '  -0.49%  '

$ws.Range('D13').Value = '1.636.95'
$ws.Range('E13').Value = '  -0.66%  '

$ws.Range('D14').Value = '4.03'
$ws.Range('E14').Value = '  -0.25%  '

$ws.Range('D15').Value = '0.557'
$ws.Range('E15').Value = '  -3.84%  '

$ws.Range('D16').Value = '64.52'
$ws.Range('E16').Value = '  +0.09%  '

$ws.Range('D17').Value = '27.611.23'
$ws.Range('E17').Value = '  +0.63%  '

$ws.Range('D18').Value = '228.62'
$ws.Range('E18').Value = '  -0.79%  '

$ws.Range('D19').Value = '7.70'
$ws.Range('E19').Value = '  +1.88%  '

$ws.Range('D22').Value = '4.29'
$ws.Range('E22').Value = '  -1.04%  '

$ws.Range('D23').Value = '10.02'
$ws.Range('E23').Value = '  +3.19%  '

$ws.Range('D24').Value = '2.00'
$ws.Range('E24').Value = '  -1.11%  '

$ws.Range('D25').Value = '150.73'
$ws.Range('E25').Value = '  +2.07%  '

$ws.Range('D26').Value = '6.95'
$ws.Range('E26').Value = '  -1.05%  '

$ws.Range('E27').Value = '  -1.96%  '

$ws.Range('E28').Value = '  +0.09%  '

$ws.Range('D29').Value = '15.56'
$ws.Range('E29').Value = '  -0.46%  '

$ws.Range('E30').Value = '  -0.24%  '

$ws.Range('E31').Value = '  -0.01%  '

$ws.Range('E32').Value = '  +0.13%  '

$ws.Range('D33').Value = '1.452.10'
$ws.Range('E33').Value = '  +2.32%  '

$ws.Range('D34').Value = '3.10'
$ws.Range('E34').Value = '  -2.15%  '

$ws.Range('E35').Value = '  -1.19%  '

$ws.Range('E36').Value = '  -0.24%  '

$ws.Range('E37').Value = '  -0.60%  '

$ws.Range('D38').Value = '0.875'
$ws.Range('E38').Value = '  -1.17%  '

$ws.Range('E39').Value = '  +0.28%  '

$ws.Range('E40').Value = '  +7.74%  '

$ws.Range('D41').Value = '69.79'
$ws.Range('E41').Value = '  +7.85%  '

$ws.Range('E42').Value = '  +0.10%  '

$ws.Range('E43').Value = '  -0.93%  '

$ws.Range('E44').Value = '  +1.43%  '

$ws.Range('B45').Value = 'mCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
$ws.Range('D45').Value = '2.47'
$ws.Range('E45').Value = '  +0.02%  '

$ws.Range('B46').Value = 'MXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D46').Value = '2.23'
$ws.Range('E46').Value = '  -0.71%  '

$ws.Range('B47').Value = 'RocketPoolETH'
$ws.Range('C47').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D47').Value = '1.778.29'
$ws.Range('E47').Value = '  -0.52%  '

$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').Value = '1.72'
$ws.Range('E48').Value = '  +2.61%  '

$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').Value = '86.30'
$ws.Range('E49').Value = '  -2.11%  '

$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '0.0₆0107'
$ws.Range('E50').Value = '  -0.31%  '

$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').Value = '0.0984'
$ws.Range('E51').Value = '  -0.97%  '
